# gaming (game_06): more assets used
# Adds 5 new columns (N:R -> sprite/spriteScale/deathSprite/deathSfx/themeTrack)
# to the bosses sheet, localizes a few existing text fields to Chinese,
# trims the "ui/assets/" prefix off the existing telegraphSprite paths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: type row for the new columns ---
$ws.Range("N4").Value = "string"
$ws.Range("O4").Value = "float"
$ws.Range("P4").Value = "string"
$ws.Range("Q4").Value = "string"
$ws.Range("R4").Value = "string"

# --- Row 5: header row for the new columns ---
$ws.Range("N5").Value = "sprite"
$ws.Range("O5").Value = "spriteScale"
$ws.Range("P5").Value = "deathSprite"
$ws.Range("Q5").Value = "deathSfx"
$ws.Range("R5").Value = "themeTrack"

# --- Row 6: The Choir of Mouths -> 万口赞歌执政体 ---
$ws.Range("D6").Value = "万口赞歌执政体"
$ws.Range("I6").Value = "旋转虚空光束配合赞歌轰炸"
$ws.Range("L6").Value = "fx/telegraph/choir_circle.png"
$ws.Range("M6").Value = "每 20 秒累积 1 层恐惧。"
$ws.Range("N6").Value = "ui/assets/topdown/top-down-shooter/characters/tank-cannon.png"
$ws.Range("O6").Value = "'1.1"
$ws.Range("O6").ClearFormats()
$ws.Range("P6").Value = "ui/assets/topdown/top-down-shooter/effects/explosion.png"
$ws.Range("Q6").Value = "ui/assets/topdown/top-down-shooter/sounds/explosion-3.wav"
$ws.Range("R6").Value = "ui/assets/topdown/top-down-shooter/music/theme-4.ogg"

# --- Row 7: Tide Shepherd -> 潮汐引航巨像 ---
$ws.Range("D7").Value = "潮汐引航巨像"
$ws.Range("I7").Value = "锁链钩拖将玩家拉入潮汐航道"
$ws.Range("L7").Value = "fx/telegraph/tidal_lane.png"
$ws.Range("M7").Value = "每 15 秒海潮横扫战场边缘。"
$ws.Range("N7").Value = "ui/assets/topdown/top-down-shooter/characters/tank-base.png"
$ws.Range("O7").Value = "'1.05"
$ws.Range("O7").ClearFormats()
$ws.Range("P7").Value = "ui/assets/topdown/top-down-shooter/effects/5.png"
$ws.Range("Q7").Value = "ui/assets/topdown/top-down-shooter/sounds/explosion-2.wav"
$ws.Range("R7").Value = "ui/assets/topdown/top-down-shooter/music/theme-4.ogg"

# --- Row 8: Nameless Beacon -> 无名肃光灯塔 ---
$ws.Range("D8").Value = "无名肃光灯塔"
$ws.Range("I8").Value = "监视立柱抽离理智能量"
$ws.Range("L8").Value = "fx/telegraph/beacon_grid.png"
$ws.Range("M8").Value = "激活时视野压缩至 65%。"
$ws.Range("N8").Value = "ui/assets/topdown/top-down-shooter/background/door.gif"
$ws.Range("O8").Value = "'1.0"
$ws.Range("O8").ClearFormats()
$ws.Range("P8").Value = "ui/assets/topdown/top-down-shooter/effects/4.png"
$ws.Range("Q8").Value = "ui/assets/topdown/top-down-shooter/sounds/explosion-1.wav"
$ws.Range("R8").Value = "ui/assets/topdown/top-down-shooter/music/theme-4.ogg"

# The whole table is authored with numbers-as-text (matches the original
# A4:M8 "ignore number stored as text" marker); tell Excel's error checker
# to keep ignoring it across the newly widened A4:R8 range too.
try {
    $ws.Range("A4:R8").Errors.Item(9).Ignore = $true
} catch {
    # Older/limited Excel automation hosts may not expose Errors.Item here;
    # the numberStoredAsText marker is cosmetic only, so ignore failures.
}
